$wb = $excel.ActiveWorkbook

# --- Sheet: compare_models (column I, "TT (Sec)") ---
$wsCompare = $wb.Worksheets.Item("compare_models")

$wsCompare.Range("I2").Value  = 0.062
$wsCompare.Range("I3").Value  = 0.044
$wsCompare.Range("I4").Value  = 0.028
$wsCompare.Range("I5").Value  = 0.084
$wsCompare.Range("I6").Value  = 0.036
$wsCompare.Range("I7").Value  = 0.018
$wsCompare.Range("I8").Value  = 1.072
$wsCompare.Range("I9").Value  = 0.022
$wsCompare.Range("I10").Value = 0.018
$wsCompare.Range("I11").Value = 0.02
$wsCompare.Range("I12").Value = 0.024
$wsCompare.Range("I13").Value = 0.504
$wsCompare.Range("I14").Value = 0.018
$wsCompare.Range("I15").Value = 0.016
$wsCompare.Range("I16").Value = 0.016
$wsCompare.Range("I17").Value = 0.016
$wsCompare.Range("I18").Value = 0.016

# --- Sheet: pred_final (row 2, columns C-H) ---
$wsPredFinal = $wb.Worksheets.Item("pred_final")

$wsPredFinal.Range("C2").Value = 1.3609
$wsPredFinal.Range("D2").Value = 4.1832
$wsPredFinal.Range("E2").Value = 2.0453
$wsPredFinal.Range("F2").Value = 0.9937
$wsPredFinal.Range("G2").Value = 0.0376
$wsPredFinal.Range("H2").Value = 0.0236
